$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 250 - this shifts existing rows 250:314
# down to 251:315 (formats copied from the row above, matching Excel's
# default Insert() behaviour).
$ws.Rows.Item(250).Insert()

# Populate the newly-inserted row 250 with the new weekly data point.
$ws.Range("A250").Value = 5
$ws.Range("B250").Value = "Macroferia Regional de Talca"
$ws.Range("C250").Value = "Maule"
$ws.Range("D250").Value = 44736
$ws.Range("E250").Value = 7
$ws.Range("F250").Value = 100112003
$ws.Range("G250").Value = "Ajo"
$ws.Range("H250").Value = "Chino"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 150
$ws.Range("K250").Value = 22000
$ws.Range("L250").Value = 22000
$ws.Range("M250").Value = 22000
$ws.Range("N250").Value = '$/malla 10 kilos'
$ws.Range("O250").Value = "China"
$ws.Range("P250").Value = 2200
$ws.Range("Q250").Value = 10
$ws.Range("R250").Value = "Hortaliza"
